# Applies the "Sample Parts-List" header-row documentation rewrite:
#  - Rewrites the 5 header-row cells (A1:E1) as bold-title + bullet-point
#    instructions (rich text: bold run for the title, regular run for body).
#  - Grows row 1's height and column B's width to fit the longer text.
#  - Left-aligns (instead of center-aligns) the header cell styles.
#  - Re-points the frozen pane / active selection like the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RichHeader([string]$CellAddr, [string]$Title, [string]$Body) {
    $cell = $ws.Range($CellAddr)
    $full = $Title + $Body
    $cell.Value2 = $full
    $cell.Characters(1, $Title.Length).Font.Bold = $true
    $cell.Characters($Title.Length + 1, $Body.Length).Font.Bold = $false
}

# --- A1: Part Number ---------------------------------------------------
$a1Body = @'

- Value required.
- Alpha-numeric up to 255 characters.
- Must be column A.
- Spaces will be removed if present.
'@
Set-RichHeader "A1" "Part Number" $a1Body

# --- B1: Reference Designators ------------------------------------------
$b1Body = @'

- Value required.
- Alpha-numeric up to 40,000 characters.
- Must be column B.
- Implicit ranges like C10-13 or C10-C13 permitted.
- Explicit ranges like C10,C11,C12,C13 permitted.
- Mixed ranges like C10-13,C17,C18 or C10-C13,C17,C18 permitted.
- Duplicate Reference Designators not permitted.
- Simple Designators like CR1, RN10, UART100, etc. permitted; letter(s) followed by number(s).
- Complex Designators like R1C2, X10Y20Z30, etc. not permitted; letter(s) followed by number(s) followed by letter(s) followed by number(s)....
- Designators containing embedded hyphens (-) or commas (,) not permitted; both are designator separators.
   - If necessary, replace with underscores (_) and semicolons (;) before processing, then restore to hyphens and commas after processing.
- Spaces will be removed if present.
'@
Set-RichHeader "B1" "Reference Designators" $b1Body

# --- C1: Part Description -----------------------------------------------
$c1Body = @'

- Value optional.
- Alpha-numeric up to 255 characters.
- Must be column C if present.
'@
Set-RichHeader "C1" "Part Description" $c1Body

# --- D1: Quantity ---------------------------------------------------------
$d1Body = @'

- Value optional.
- Numeric integer in inclusive range {1..32,767}.
- Must be column D if present.
- Used to optionally compare stated Quantity to calculated Quantity.
'@
Set-RichHeader "D1" "Quantity" $d1Body

# --- E1: Comments -----------------------------------------------------------
$e1Body = @'

- Parts-List Processor ignores the 1st row during processing, so a place-holding header row is required.
- This row is the default header row and processes fine, but can be customized as desired.
- A blank 1st row processes fine also.
'@
Set-RichHeader "E1" "Comments" $e1Body

# --- Row / column sizing -------------------------------------------------
$ws.Rows(1).RowHeight = 220.8
$ws.Columns(2).ColumnWidth = 92.88671875

# --- Header cell alignment: center -> left --------------------------------
$ws.Range("A1:E1").HorizontalAlignment = -4131  # xlLeft

# --- View: frozen pane anchor + active selection --------------------------
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 18
